$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update subject id header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2: the CON value moved from D2 to C2
$ws.Range("D2").ClearContents()
$ws.Range("C2").Value = 30.373805491377226

# Row 3: clear out the old B3/C3 values (no replacement in this row)
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()

# Update the selection to match the new highlighted range
$ws.Range("B1:E3").Select()
